# Prep labels complete without indivisible pack error
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that are "Already Completed" (0 packs, not actionable) or otherwise
# dropped from the final prep list. Delete from bottom to top so row
# numbers of not-yet-deleted rows stay stable.
$rowsToDelete = @(31, 11, 10, 8, 6, 5, 4, 3, 2)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# After the deletions, 21 data rows remain (rows 2-22). Four of them
# encountered an indivisible-pack error during label prep and need their
# Status updated from "Complete" to "Error" (one of them also had its
# pack count reset to 0 since the packs could not be processed).
$ws.Range("E3").Value = "Error"   # CMH2 35247212391
$ws.Range("E8").Value = "Error"   # OAK3 35247173221
$ws.Range("D11").Value = 0        # PHX7 35247173201 pack count -> 0
$ws.Range("E11").Value = "Error"  # PHX7 35247173201
$ws.Range("E17").Value = "Error"  # SNA4 35247212451

$ws.Range("A1").Select()
